$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'29.533.02"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'  +0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'1.915.01"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'  -0.18%  "
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'  +0.72%  "
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'326.27"
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'  +0.10%  "
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'  +0.63%  "
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'  +1.68%  "
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'0.4072"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'  -0.56%  "
$c.Style = 'Normal'
$c = $ws.Range('B9')
$c.Value = "'Dogecoin"
$c.Style = 'Normal'
$c = $ws.Range('C9')
$c.Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'0.08144"
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'  +1.14%  "
$c.Style = 'Normal'
$c = $ws.Range('B10')
$c.Value = "'Polygon"
$c.Style = 'Normal'
$c = $ws.Range('C10')
$c.Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'1.013"
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'  +0.09%  "
$c.Style = 'Normal'
$c = $ws.Range('B11')
$c.Value = "'Solana"
$c.Style = 'Normal'
$c = $ws.Range('C11')
$c.Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'23.39"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'  +3.96%  "
$c.Style = 'Normal'
$c = $ws.Range('B12')
$c.Value = "'WrappedEther"
$c.Style = 'Normal'
$c = $ws.Range('C12')
$c.Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'1.933.95"
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'  +1.56%  "
$c.Style = 'Normal'
$c = $ws.Range('B13')
$c.Value = "'Polkadot"
$c.Style = 'Normal'
$c = $ws.Range('C13')
$c.Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'5.998"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'  +1.09%  "
$c.Style = 'Normal'
$c = $ws.Range('B14')
$c.Value = "'Chainlink"
$c.Style = 'Normal'
$c = $ws.Range('C14')
$c.Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'7.130"
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'  -0.41%  "
$c.Style = 'Normal'
$c = $ws.Range('B15')
$c.Value = "'Litecoin"
$c.Style = 'Normal'
$c = $ws.Range('C15')
$c.Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'90.26"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'  +0.73%  "
$c.Style = 'Normal'
$c = $ws.Range('B16')
$c.Value = "'TRON"
$c.Style = 'Normal'
$c = $ws.Range('C16')
$c.Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'0.06787"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'  +2.86%  "
$c.Style = 'Normal'
$c = $ws.Range('B17')
$c.Value = "'BinanceUSD"
$c.Style = 'Normal'
$c = $ws.Range('C17')
$c.Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'1.008"
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'  +0.68%  "
$c.Style = 'Normal'
$c = $ws.Range('B18')
$c.Value = "'ShibaInu"
$c.Style = 'Normal'
$c = $ws.Range('C18')
$c.Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'0.00001039"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'  +0.58%  "
$c.Style = 'Normal'
$c = $ws.Range('B19')
$c.Value = "'Avalanche"
$c.Style = 'Normal'
$c = $ws.Range('C19')
$c.Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = "'17.70"
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'  -0.56%  "
$c.Style = 'Normal'
$c = $ws.Range('B20')
$c.Value = "'Dai"
$c.Style = 'Normal'
$c = $ws.Range('C20')
$c.Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = "'1.007"
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'  +0.72%  "
$c.Style = 'Normal'
$c = $ws.Range('B21')
$c.Value = "'WrappedBTC"
$c.Style = 'Normal'
$c = $ws.Range('C21')
$c.Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.Value = "'29.541.41"
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'  +0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('B22')
$c.Value = "'Uniswap"
$c.Style = 'Normal'
$c = $ws.Range('C22')
$c.Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'5.620"
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'  +1.32%  "
$c.Style = 'Normal'
$c = $ws.Range('B23')
$c.Value = "'Cosmos"
$c.Style = 'Normal'
$c = $ws.Range('C23')
$c.Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'11.79"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'  +2.56%  "
$c.Style = 'Normal'
$c = $ws.Range('B24')
$c.Value = "'Toncoin"
$c.Style = 'Normal'
$c = $ws.Range('C24')
$c.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = "'2.182"
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'  -1.16%  "
$c.Style = 'Normal'
$c = $ws.Range('B25')
$c.Value = "'WrappedliquidstakedEther2.0"
$c.Style = 'Normal'
$c = $ws.Range('C25')
$c.Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'2.151.26"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'  +0.86%  "
$c.Style = 'Normal'
$c = $ws.Range('B26')
$c.Value = "'Monero"
$c.Style = 'Normal'
$c = $ws.Range('C26')
$c.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = "'155.30"
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'  +0.39%  "
$c.Style = 'Normal'
$c = $ws.Range('B27')
$c.Value = "'InternetComputer(DFINITY)"
$c.Style = 'Normal'
$c = $ws.Range('C27')
$c.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'6.487"
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = "'  +6.99%  "
$c.Style = 'Normal'
$c = $ws.Range('B28')
$c.Value = "'EthereumClassic"
$c.Style = 'Normal'
$c = $ws.Range('C28')
$c.Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = "'20.07"
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = "'  +1.14%  "
$c.Style = 'Normal'
$c = $ws.Range('B29')
$c.Value = "'LidoDAOToken"
$c.Style = 'Normal'
$c = $ws.Range('C29')
$c.Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.Value = "'2.099"
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.Value = "'  -1.49%  "
$c.Style = 'Normal'
$c = $ws.Range('B30')
$c.Value = "'BitcoinCash"
$c.Style = 'Normal'
$c = $ws.Range('C30')
$c.Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.Value = "'119.68"
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.Value = "'  +1.50%  "
$c.Style = 'Normal'
$c = $ws.Range('B31')
$c.Value = "'ImmutableX"
$c.Style = 'Normal'
$c = $ws.Range('C31')
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'1.032"
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = "'  -2.95%  "
$c.Style = 'Normal'
$c = $ws.Range('B32')
$c.Value = "'Stellar"
$c.Style = 'Normal'
$c = $ws.Range('C32')
$c.Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = "'0.09557"
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = "'  +0.19%  "
$c.Style = 'Normal'
$c = $ws.Range('B33')
$c.Value = "'Filecoin"
$c.Style = 'Normal'
$c = $ws.Range('C33')
$c.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.Value = "'5.519"
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = "'  +2.00%  "
$c.Style = 'Normal'
$c = $ws.Range('B34')
$c.Value = "'HuobiToken"
$c.Style = 'Normal'
$c = $ws.Range('C34')
$c.Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.Value = "'3.566"
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = "'  +0.23%  "
$c.Style = 'Normal'
$c = $ws.Range('B35')
$c.Value = "'ARBITRUM"
$c.Style = 'Normal'
$c = $ws.Range('C35')
$c.Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.Value = "'1.395"
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = "'  -2.75%  "
$c.Style = 'Normal'
$c = $ws.Range('B36')
$c.Value = "'VeChain"
$c.Style = 'Normal'
$c = $ws.Range('C36')
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.Value = "'0.02269"
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.Value = "'  +0.41%  "
$c.Style = 'Normal'
$c = $ws.Range('B37')
$c.Value = "'Hedera"
$c.Style = 'Normal'
$c = $ws.Range('C37')
$c.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.Value = "'0.06100"
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = "'  -0.36%  "
$c.Style = 'Normal'
$c = $ws.Range('B38')
$c.Value = "'TrustWalletToken"
$c.Style = 'Normal'
$c = $ws.Range('C38')
$c.Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = "'1.182"
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'  +0.90%  "
$c.Style = 'Normal'
$c = $ws.Range('B39')
$c.Value = "'TheSandbox"
$c.Style = 'Normal'
$c = $ws.Range('C39')
$c.Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'0.5939"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'  +0.71%  "
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'  +6.07%  "
$c.Style = 'Normal'
$c = $ws.Range('B41')
$c.Value = "'FraxShare"
$c.Style = 'Normal'
$c = $ws.Range('C41')
$c.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'7.969"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'  -4.31%  "
$c.Style = 'Normal'
$c = $ws.Range('B42')
$c.Value = "'Algorand"
$c.Style = 'Normal'
$c = $ws.Range('C42')
$c.Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.Value = "'0.1857"
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'  +0.60%  "
$c.Style = 'Normal'
$c = $ws.Range('B43')
$c.Value = "'RenderToken"
$c.Style = 'Normal'
$c = $ws.Range('C43')
$c.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.Value = "'2.480"
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'  -3.17%  "
$c.Style = 'Normal'
$c = $ws.Range('B44')
$c.Value = "'Cronos"
$c.Style = 'Normal'
$c = $ws.Range('C44')
$c.Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.Value = "'0.07718"
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'  -3.78%  "
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.Value = "'1.248"
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'  -2.92%  "
$c.Style = 'Normal'
$c = $ws.Range('B46')
$c.Value = "'EnergySwap"
$c.Style = 'Normal'
$c = $ws.Range('C46')
$c.Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.Value = "'12.40"
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'  +2.26%  "
$c.Style = 'Normal'
$c = $ws.Range('B47')
$c.Value = "'Decentraland"
$c.Style = 'Normal'
$c = $ws.Range('C47')
$c.Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'0.5570"
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = "'  +0.12%  "
$c.Style = 'Normal'
$c = $ws.Range('B48')
$c.Value = "'NEARProtocol"
$c.Style = 'Normal'
$c = $ws.Range('C48')
$c.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'1.943"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'  +0.20%  "
$c.Style = 'Normal'
$c = $ws.Range('B49')
$c.Value = "'Quant"
$c.Style = 'Normal'
$c = $ws.Range('C49')
$c.Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'115.67"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'  +2.22%  "
$c.Style = 'Normal'
$c = $ws.Range('B50')
$c.Value = "'Aave"
$c.Style = 'Normal'
$c = $ws.Range('C50')
$c.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'72.76"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'  +1.47%  "
$c.Style = 'Normal'
$c = $ws.Range('B51')
$c.Value = "'EOS"
$c.Style = 'Normal'
$c = $ws.Range('C51')
$c.Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = "'1.053"
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'  +1.87%  "
$c.Style = 'Normal'
